# Update division exercises in the table to new values
$d = $word.ActiveDocument

$replacements = @(
    @("287÷8=", "818÷7="),
    @("713÷9=", "440÷5="),
    @("644÷2=", "737÷3="),
    @("122÷3=", "545÷5="),
    @("387÷9=", "358÷6="),
    @("432÷7=", "427÷3="),
    @("464÷6=", "104÷7="),
    @("513÷3=", "402÷5="),
    @("156÷9=", "557÷8="),
    @("425÷4=", "604÷8="),
    @("684÷9=", "107÷2="),
    @("749÷2=", "482÷5="),
    @("228÷2=", "540÷3="),
    @("996÷9=", "382÷6="),
    @("371÷9=", "628÷5="),
    @("602÷3=", "846÷3="),
    @("522÷5=", "978÷2="),
    @("127÷2=", "183÷8="),
    @("744÷5=", "198÷2="),
    @("769÷5=", "364÷7="),
    @("120÷3=", "824÷7="),
    @("984÷2=", "938÷5="),
    @("560÷5=", "502÷4="),
    @("336÷2=", "609÷9="),
    @("582÷6=", "917÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
